# Auto-generated Excel COM-interop script to apply the betexplorer data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: re-shuffle rows whose match data moved to a different row of the same match-day block ---

# Row 34: data now matches former row 36
$ws.Cells.Item(34, 6).Value = "Chelsea"
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = "Nottingham"
$ws.Cells.Item(34, 9).Value = 1
$ws.Cells.Item(34, 10).Value = 1.51
$ws.Cells.Item(34, 11).Value = "20/08/2023 09:02"
$ws.Cells.Item(34, 12).Value = 1.41
$ws.Cells.Item(34, 13).Value = "02/09/2023 15:50"
$ws.Cells.Item(34, 14).Value = 4.47
$ws.Cells.Item(34, 15).Value = "20/08/2023 09:02"
$ws.Cells.Item(34, 16).Value = 5.01
$ws.Cells.Item(34, 17).Value = "02/09/2023 15:58"
$ws.Cells.Item(34, 18).Value = 6.71
$ws.Cells.Item(34, 19).Value = "20/08/2023 09:02"
$ws.Cells.Item(34, 20).Value = 8.449999999999999
$ws.Cells.Item(34, 21).Value = "02/09/2023 15:58"
$ws.Cells.Item(34, 22).Value = "https://www.betexplorer.com/football/england/premier-league/chelsea-nottingham/0d8k37tt/"

# Row 35: data now matches former row 34
$ws.Cells.Item(35, 6).Value = "Burnley"
$ws.Cells.Item(35, 7).Value = 2
$ws.Cells.Item(35, 8).Value = "Tottenham"
$ws.Cells.Item(35, 9).Value = 5
$ws.Cells.Item(35, 10).Value = 3.81
$ws.Cells.Item(35, 11).Value = "21/08/2023 06:14"
$ws.Cells.Item(35, 12).Value = 3.71
$ws.Cells.Item(35, 13).Value = "02/09/2023 15:53"
$ws.Cells.Item(35, 14).Value = 3.57
$ws.Cells.Item(35, 15).Value = "21/08/2023 06:14"
$ws.Cells.Item(35, 16).Value = 3.79
$ws.Cells.Item(35, 17).Value = "02/09/2023 15:59"
$ws.Cells.Item(35, 18).Value = 1.94
$ws.Cells.Item(35, 19).Value = "21/08/2023 06:14"
$ws.Cells.Item(35, 20).Value = 2.04
$ws.Cells.Item(35, 21).Value = "02/09/2023 15:52"
$ws.Cells.Item(35, 22).Value = "https://www.betexplorer.com/football/england/premier-league/burnley-tottenham/E7jreAlJ/"

# Row 36: data now matches former row 35
$ws.Cells.Item(36, 6).Value = "Brentford"
$ws.Cells.Item(36, 7).Value = 2
$ws.Cells.Item(36, 8).Value = "Bournemouth"
$ws.Cells.Item(36, 9).Value = 2
$ws.Cells.Item(36, 10).Value = 1.73
$ws.Cells.Item(36, 11).Value = "20/08/2023 09:02"
$ws.Cells.Item(36, 12).Value = 1.75
$ws.Cells.Item(36, 13).Value = "02/09/2023 15:27"
$ws.Cells.Item(36, 14).Value = 4.37
$ws.Cells.Item(36, 15).Value = "20/08/2023 09:02"
$ws.Cells.Item(36, 16).Value = 4.01
$ws.Cells.Item(36, 17).Value = "02/09/2023 15:30"
$ws.Cells.Item(36, 18).Value = 3.95
$ws.Cells.Item(36, 19).Value = "20/08/2023 09:02"
$ws.Cells.Item(36, 20).Value = 4.86
$ws.Cells.Item(36, 21).Value = "02/09/2023 15:59"
$ws.Cells.Item(36, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brentford-bournemouth/hh2ZdWJ6/"

# Row 42: data now matches former row 46
$ws.Cells.Item(42, 6).Value = "West Ham"
$ws.Cells.Item(42, 7).Value = 1
$ws.Cells.Item(42, 8).Value = "Manchester City"
$ws.Cells.Item(42, 9).Value = 3
$ws.Cells.Item(42, 10).Value = 6.22
$ws.Cells.Item(42, 11).Value = "28/08/2023 09:02"
$ws.Cells.Item(42, 12).Value = 5.49
$ws.Cells.Item(42, 13).Value = "16/09/2023 15:58"
$ws.Cells.Item(42, 14).Value = 4.81
$ws.Cells.Item(42, 15).Value = "28/08/2023 09:02"
$ws.Cells.Item(42, 16).Value = 4.3
$ws.Cells.Item(42, 17).Value = "16/09/2023 15:59"
$ws.Cells.Item(42, 18).Value = 1.44
$ws.Cells.Item(42, 19).Value = "28/08/2023 09:02"
$ws.Cells.Item(42, 20).Value = 1.63
$ws.Cells.Item(42, 21).Value = "16/09/2023 15:54"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/england/premier-league/west-ham-manchester-city/UcysC2PA/"

# Row 43: data now matches former row 45
$ws.Cells.Item(43, 6).Value = "Tottenham"
$ws.Cells.Item(43, 7).Value = 2
$ws.Cells.Item(43, 8).Value = "Sheffield Utd"
$ws.Cells.Item(43, 9).Value = 1
$ws.Cells.Item(43, 10).Value = 1.29
$ws.Cells.Item(43, 11).Value = "28/08/2023 11:55"
$ws.Cells.Item(43, 12).Value = 1.3
$ws.Cells.Item(43, 13).Value = "16/09/2023 15:58"
$ws.Cells.Item(43, 14).Value = 5.56
$ws.Cells.Item(43, 15).Value = "28/08/2023 11:55"
$ws.Cells.Item(43, 16).Value = 6.31
$ws.Cells.Item(43, 17).Value = "16/09/2023 15:58"
$ws.Cells.Item(43, 18).Value = 9.48
$ws.Cells.Item(43, 19).Value = "28/08/2023 11:55"
$ws.Cells.Item(43, 20).Value = 9.4
$ws.Cells.Item(43, 21).Value = "16/09/2023 15:58"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/england/premier-league/tottenham-sheffield-utd/nZoxDrA4/"

# Row 44: data now matches former row 43
$ws.Cells.Item(44, 6).Value = "Fulham"
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = "Luton"
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 1.68
$ws.Cells.Item(44, 11).Value = "28/08/2023 11:54"
$ws.Cells.Item(44, 12).Value = 1.69
$ws.Cells.Item(44, 13).Value = "16/09/2023 15:36"
$ws.Cells.Item(44, 14).Value = 3.8
$ws.Cells.Item(44, 15).Value = "28/08/2023 11:54"
$ws.Cells.Item(44, 16).Value = 3.82
$ws.Cells.Item(44, 17).Value = "16/09/2023 15:59"
$ws.Cells.Item(44, 18).Value = 4.97
$ws.Cells.Item(44, 19).Value = "28/08/2023 11:54"
$ws.Cells.Item(44, 20).Value = 5.68
$ws.Cells.Item(44, 21).Value = "16/09/2023 15:59"
$ws.Cells.Item(44, 22).Value = "https://www.betexplorer.com/football/england/premier-league/fulham-luton/bD5si1mo/"

# Row 45: data now matches former row 42
$ws.Cells.Item(45, 6).Value = "Aston Villa"
$ws.Cells.Item(45, 7).Value = 3
$ws.Cells.Item(45, 8).Value = "Crystal Palace"
$ws.Cells.Item(45, 9).Value = 1
$ws.Cells.Item(45, 10).Value = 1.88
$ws.Cells.Item(45, 11).Value = "28/08/2023 09:02"
$ws.Cells.Item(45, 12).Value = 1.98
$ws.Cells.Item(45, 13).Value = "16/09/2023 15:53"
$ws.Cells.Item(45, 14).Value = 3.72
$ws.Cells.Item(45, 15).Value = "28/08/2023 09:02"
$ws.Cells.Item(45, 16).Value = 3.66
$ws.Cells.Item(45, 17).Value = "16/09/2023 15:53"
$ws.Cells.Item(45, 18).Value = 3.91
$ws.Cells.Item(45, 19).Value = "28/08/2023 09:02"
$ws.Cells.Item(45, 20).Value = 4.08
$ws.Cells.Item(45, 21).Value = "16/09/2023 15:53"
$ws.Cells.Item(45, 22).Value = "https://www.betexplorer.com/football/england/premier-league/aston-villa-crystal-palace/SUEBdNPN/"

# Row 46: data now matches former row 44
$ws.Cells.Item(46, 6).Value = "Manchester Utd"
$ws.Cells.Item(46, 7).Value = 1
$ws.Cells.Item(46, 8).Value = "Brighton"
$ws.Cells.Item(46, 9).Value = 3
$ws.Cells.Item(46, 10).Value = 1.88
$ws.Cells.Item(46, 11).Value = "28/08/2023 09:02"
$ws.Cells.Item(46, 12).Value = 2.07
$ws.Cells.Item(46, 13).Value = "16/09/2023 15:59"
$ws.Cells.Item(46, 14).Value = 3.99
$ws.Cells.Item(46, 15).Value = "28/08/2023 09:02"
$ws.Cells.Item(46, 16).Value = 3.9
$ws.Cells.Item(46, 17).Value = "16/09/2023 15:59"
$ws.Cells.Item(46, 18).Value = 3.6
$ws.Cells.Item(46, 19).Value = "28/08/2023 09:02"
$ws.Cells.Item(46, 20).Value = 3.53
$ws.Cells.Item(46, 21).Value = "16/09/2023 15:59"
$ws.Cells.Item(46, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-united-brighton/0IqQFpQo/"

# Row 51: data now matches former row 53
$ws.Cells.Item(51, 6).Value = "Manchester City"
$ws.Cells.Item(51, 7).Value = 2
$ws.Cells.Item(51, 8).Value = "Nottingham"
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 1.12
$ws.Cells.Item(51, 11).Value = "05/09/2023 12:01"
$ws.Cells.Item(51, 12).Value = 1.17
$ws.Cells.Item(51, 13).Value = "23/09/2023 15:28"
$ws.Cells.Item(51, 14).Value = 9.08
$ws.Cells.Item(51, 15).Value = "05/09/2023 12:01"
$ws.Cells.Item(51, 16).Value = 8.5
$ws.Cells.Item(51, 17).Value = "23/09/2023 15:17"
$ws.Cells.Item(51, 18).Value = 15.94
$ws.Cells.Item(51, 19).Value = "05/09/2023 12:01"
$ws.Cells.Item(51, 20).Value = 18
$ws.Cells.Item(51, 21).Value = "23/09/2023 15:28"
$ws.Cells.Item(51, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-city-nottingham/GSENOu9G/"

# Row 53: data now matches former row 51
$ws.Cells.Item(53, 6).Value = "Crystal Palace"
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = "Fulham"
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 1.88
$ws.Cells.Item(53, 11).Value = "05/09/2023 12:01"
$ws.Cells.Item(53, 12).Value = 2.1
$ws.Cells.Item(53, 13).Value = "23/09/2023 15:59"
$ws.Cells.Item(53, 14).Value = 3.55
$ws.Cells.Item(53, 15).Value = "05/09/2023 12:01"
$ws.Cells.Item(53, 16).Value = 3.39
$ws.Cells.Item(53, 17).Value = "23/09/2023 15:59"
$ws.Cells.Item(53, 18).Value = 4.11
$ws.Cells.Item(53, 19).Value = "05/09/2023 12:01"
$ws.Cells.Item(53, 20).Value = 3.97
$ws.Cells.Item(53, 21).Value = "23/09/2023 15:59"
$ws.Cells.Item(53, 22).Value = "https://www.betexplorer.com/football/england/premier-league/crystal-palace-fulham/0xCBRsPc/"

# Row 57: data now matches former row 59
$ws.Cells.Item(57, 6).Value = "Liverpool"
$ws.Cells.Item(57, 7).Value = 3
$ws.Cells.Item(57, 8).Value = "West Ham"
$ws.Cells.Item(57, 9).Value = 1
$ws.Cells.Item(57, 10).Value = 1.31
$ws.Cells.Item(57, 11).Value = "05/09/2023 12:01"
$ws.Cells.Item(57, 12).Value = 1.39
$ws.Cells.Item(57, 13).Value = "24/09/2023 14:53"
$ws.Cells.Item(57, 14).Value = 5.93
$ws.Cells.Item(57, 15).Value = "05/09/2023 12:01"
$ws.Cells.Item(57, 16).Value = 5.64
$ws.Cells.Item(57, 17).Value = "24/09/2023 14:58"
$ws.Cells.Item(57, 18).Value = 7.64
$ws.Cells.Item(57, 19).Value = "05/09/2023 12:01"
$ws.Cells.Item(57, 20).Value = 7.61
$ws.Cells.Item(57, 21).Value = "24/09/2023 14:59"
$ws.Cells.Item(57, 22).Value = "https://www.betexplorer.com/football/england/premier-league/liverpool-west-ham/r11GQ1v4/"

# Row 59: data now matches former row 57
$ws.Cells.Item(59, 6).Value = "Brighton"
$ws.Cells.Item(59, 7).Value = 3
$ws.Cells.Item(59, 8).Value = "Bournemouth"
$ws.Cells.Item(59, 9).Value = 1
$ws.Cells.Item(59, 10).Value = 1.32
$ws.Cells.Item(59, 11).Value = "05/09/2023 12:01"
$ws.Cells.Item(59, 12).Value = 1.52
$ws.Cells.Item(59, 13).Value = "24/09/2023 14:30"
$ws.Cells.Item(59, 14).Value = 6.03
$ws.Cells.Item(59, 15).Value = "05/09/2023 12:01"
$ws.Cells.Item(59, 16).Value = 4.97
$ws.Cells.Item(59, 17).Value = "24/09/2023 14:53"
$ws.Cells.Item(59, 18).Value = 8.84
$ws.Cells.Item(59, 19).Value = "05/09/2023 12:01"
$ws.Cells.Item(59, 20).Value = 5.9
$ws.Cells.Item(59, 21).Value = "24/09/2023 14:59"
$ws.Cells.Item(59, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brighton-bournemouth/fuLL4KHp/"

# Row 62: data now matches former row 65
$ws.Cells.Item(62, 6).Value = "Manchester Utd"
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = "Crystal Palace"
$ws.Cells.Item(62, 9).Value = 1
$ws.Cells.Item(62, 10).Value = 1.5
$ws.Cells.Item(62, 11).Value = "17/09/2023 09:01"
$ws.Cells.Item(62, 12).Value = 1.58
$ws.Cells.Item(62, 13).Value = "30/09/2023 15:58"
$ws.Cells.Item(62, 14).Value = 4.63
$ws.Cells.Item(62, 15).Value = "17/09/2023 09:01"
$ws.Cells.Item(62, 16).Value = 4.37
$ws.Cells.Item(62, 17).Value = "30/09/2023 15:59"
$ws.Cells.Item(62, 18).Value = 6.58
$ws.Cells.Item(62, 19).Value = "17/09/2023 09:01"
$ws.Cells.Item(62, 20).Value = 5.96
$ws.Cells.Item(62, 21).Value = "30/09/2023 15:59"
$ws.Cells.Item(62, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-united-crystal-palace/Q1OnP9Kd/"

# Row 64: data now matches former row 62
$ws.Cells.Item(64, 6).Value = "West Ham"
$ws.Cells.Item(64, 7).Value = 2
$ws.Cells.Item(64, 8).Value = "Sheffield Utd"
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 1.49
$ws.Cells.Item(64, 11).Value = "23/09/2023 17:42"
$ws.Cells.Item(64, 12).Value = 1.48
$ws.Cells.Item(64, 13).Value = "30/09/2023 15:46"
$ws.Cells.Item(64, 14).Value = 4.43
$ws.Cells.Item(64, 15).Value = "23/09/2023 17:42"
$ws.Cells.Item(64, 16).Value = 4.79
$ws.Cells.Item(64, 17).Value = "30/09/2023 15:55"
$ws.Cells.Item(64, 18).Value = 5.97
$ws.Cells.Item(64, 19).Value = "23/09/2023 17:42"
$ws.Cells.Item(64, 20).Value = 7
$ws.Cells.Item(64, 21).Value = "30/09/2023 15:59"
$ws.Cells.Item(64, 22).Value = "https://www.betexplorer.com/football/england/premier-league/west-ham-sheffield-utd/hEF3LRJL/"

# Row 65: data now matches former row 66
$ws.Cells.Item(65, 6).Value = "Bournemouth"
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = "Arsenal"
$ws.Cells.Item(65, 9).Value = 4
$ws.Cells.Item(65, 10).Value = 6.39
$ws.Cells.Item(65, 11).Value = "17/09/2023 09:01"
$ws.Cells.Item(65, 12).Value = 6.06
$ws.Cells.Item(65, 13).Value = "30/09/2023 15:58"
$ws.Cells.Item(65, 14).Value = 4.83
$ws.Cells.Item(65, 15).Value = "17/09/2023 09:01"
$ws.Cells.Item(65, 16).Value = 4.62
$ws.Cells.Item(65, 17).Value = "30/09/2023 15:59"
$ws.Cells.Item(65, 18).Value = 1.49
$ws.Cells.Item(65, 19).Value = "17/09/2023 09:01"
$ws.Cells.Item(65, 20).Value = 1.55
$ws.Cells.Item(65, 21).Value = "30/09/2023 15:53"
$ws.Cells.Item(65, 22).Value = "https://www.betexplorer.com/football/england/premier-league/bournemouth-arsenal/xr3WMJwT/"

# Row 66: data now matches former row 67
$ws.Cells.Item(66, 6).Value = "Everton"
$ws.Cells.Item(66, 7).Value = 1
$ws.Cells.Item(66, 8).Value = "Luton"
$ws.Cells.Item(66, 9).Value = 2
$ws.Cells.Item(66, 10).Value = 1.58
$ws.Cells.Item(66, 11).Value = "23/09/2023 17:43"
$ws.Cells.Item(66, 12).Value = 1.67
$ws.Cells.Item(66, 13).Value = "30/09/2023 15:59"
$ws.Cells.Item(66, 14).Value = 3.94
$ws.Cells.Item(66, 15).Value = "23/09/2023 17:43"
$ws.Cells.Item(66, 16).Value = 3.93
$ws.Cells.Item(66, 17).Value = "30/09/2023 15:54"
$ws.Cells.Item(66, 18).Value = 5.63
$ws.Cells.Item(66, 19).Value = "23/09/2023 17:43"
$ws.Cells.Item(66, 20).Value = 5.7
$ws.Cells.Item(66, 21).Value = "30/09/2023 15:59"
$ws.Cells.Item(66, 22).Value = "https://www.betexplorer.com/football/england/premier-league/everton-luton/8Qg2Hc1j/"

# Row 67: data now matches former row 64
$ws.Cells.Item(67, 6).Value = "Newcastle"
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = "Burnley"
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 1.34
$ws.Cells.Item(67, 11).Value = "23/09/2023 17:43"
$ws.Cells.Item(67, 12).Value = 1.38
$ws.Cells.Item(67, 13).Value = "30/09/2023 15:56"
$ws.Cells.Item(67, 14).Value = 5.32
$ws.Cells.Item(67, 15).Value = "23/09/2023 17:43"
$ws.Cells.Item(67, 16).Value = 5.24
$ws.Cells.Item(67, 17).Value = "30/09/2023 15:45"
$ws.Cells.Item(67, 18).Value = 7.68
$ws.Cells.Item(67, 19).Value = "23/09/2023 17:43"
$ws.Cells.Item(67, 20).Value = 8.529999999999999
$ws.Cells.Item(67, 21).Value = "30/09/2023 15:58"
$ws.Cells.Item(67, 22).Value = "https://www.betexplorer.com/football/england/premier-league/newcastle-utd-burnley/4hNjOTZ2/"

# Row 73: data now matches former row 74
$ws.Cells.Item(73, 6).Value = "Burnley"
$ws.Cells.Item(73, 7).Value = 1
$ws.Cells.Item(73, 8).Value = "Chelsea"
$ws.Cells.Item(73, 9).Value = 4
$ws.Cells.Item(73, 10).Value = 4.24
$ws.Cells.Item(73, 11).Value = "28/09/2023 14:25"
$ws.Cells.Item(73, 12).Value = 5.2
$ws.Cells.Item(73, 13).Value = "07/10/2023 16:00"
$ws.Cells.Item(73, 14).Value = 3.78
$ws.Cells.Item(73, 15).Value = "28/09/2023 14:25"
$ws.Cells.Item(73, 16).Value = 4.07
$ws.Cells.Item(73, 17).Value = "07/10/2023 15:58"
$ws.Cells.Item(73, 18).Value = 1.79
$ws.Cells.Item(73, 19).Value = "28/09/2023 14:25"
$ws.Cells.Item(73, 20).Value = 1.69
$ws.Cells.Item(73, 21).Value = "07/10/2023 15:58"
$ws.Cells.Item(73, 22).Value = "https://www.betexplorer.com/football/england/premier-league/burnley-chelsea/pCfrEqCe/"

# Row 74: data now matches former row 75
$ws.Cells.Item(74, 6).Value = "Everton"
$ws.Cells.Item(74, 7).Value = 3
$ws.Cells.Item(74, 8).Value = "Bournemouth"
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 1.69
$ws.Cells.Item(74, 11).Value = "24/09/2023 10:02"
$ws.Cells.Item(74, 12).Value = 1.93
$ws.Cells.Item(74, 13).Value = "07/10/2023 15:58"
$ws.Cells.Item(74, 14).Value = 4.01
$ws.Cells.Item(74, 15).Value = "24/09/2023 10:02"
$ws.Cells.Item(74, 16).Value = 3.78
$ws.Cells.Item(74, 17).Value = "07/10/2023 15:58"
$ws.Cells.Item(74, 18).Value = 4.57
$ws.Cells.Item(74, 19).Value = "24/09/2023 10:02"
$ws.Cells.Item(74, 20).Value = 4.18
$ws.Cells.Item(74, 21).Value = "07/10/2023 15:58"
$ws.Cells.Item(74, 22).Value = "https://www.betexplorer.com/football/england/premier-league/everton-bournemouth/CInUym42/"

# Row 75: data now matches former row 73
$ws.Cells.Item(75, 6).Value = "Manchester Utd"
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = "Brentford"
$ws.Cells.Item(75, 9).Value = 1
$ws.Cells.Item(75, 10).Value = 1.55
$ws.Cells.Item(75, 11).Value = "24/09/2023 10:01"
$ws.Cells.Item(75, 12).Value = 1.63
$ws.Cells.Item(75, 13).Value = "07/10/2023 15:50"
$ws.Cells.Item(75, 14).Value = 4.44
$ws.Cells.Item(75, 15).Value = "24/09/2023 10:01"
$ws.Cells.Item(75, 16).Value = 4.33
$ws.Cells.Item(75, 17).Value = "07/10/2023 15:53"
$ws.Cells.Item(75, 18).Value = 5.18
$ws.Cells.Item(75, 19).Value = "24/09/2023 10:01"
$ws.Cells.Item(75, 20).Value = 5.41
$ws.Cells.Item(75, 21).Value = "07/10/2023 15:54"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-united-brentford/8pQbIb3s/"

# Row 78: data now matches former row 80
$ws.Cells.Item(78, 6).Value = "Brighton"
$ws.Cells.Item(78, 7).Value = 2
$ws.Cells.Item(78, 8).Value = "Liverpool"
$ws.Cells.Item(78, 9).Value = 2
$ws.Cells.Item(78, 10).Value = 2.47
$ws.Cells.Item(78, 11).Value = "24/09/2023 10:01"
$ws.Cells.Item(78, 12).Value = 3.26
$ws.Cells.Item(78, 13).Value = "08/10/2023 14:44"
$ws.Cells.Item(78, 14).Value = 3.9
$ws.Cells.Item(78, 15).Value = "24/09/2023 10:01"
$ws.Cells.Item(78, 16).Value = 4.19
$ws.Cells.Item(78, 17).Value = "08/10/2023 14:44"
$ws.Cells.Item(78, 18).Value = 2.54
$ws.Cells.Item(78, 19).Value = "24/09/2023 10:01"
$ws.Cells.Item(78, 20).Value = 2.1
$ws.Cells.Item(78, 21).Value = "08/10/2023 14:44"
$ws.Cells.Item(78, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brighton-liverpool/2m5wFPdk/"

# Row 79: data now matches former row 78
$ws.Cells.Item(79, 6).Value = "West Ham"
$ws.Cells.Item(79, 7).Value = 2
$ws.Cells.Item(79, 8).Value = "Newcastle"
$ws.Cells.Item(79, 9).Value = 2
$ws.Cells.Item(79, 10).Value = 3.81
$ws.Cells.Item(79, 11).Value = "24/09/2023 10:01"
$ws.Cells.Item(79, 12).Value = 3.58
$ws.Cells.Item(79, 13).Value = "08/10/2023 14:55"
$ws.Cells.Item(79, 14).Value = 3.87
$ws.Cells.Item(79, 15).Value = "24/09/2023 10:01"
$ws.Cells.Item(79, 16).Value = 3.73
$ws.Cells.Item(79, 17).Value = "08/10/2023 14:59"
$ws.Cells.Item(79, 18).Value = 1.86
$ws.Cells.Item(79, 19).Value = "24/09/2023 10:01"
$ws.Cells.Item(79, 20).Value = 2.11
$ws.Cells.Item(79, 21).Value = "08/10/2023 14:59"
$ws.Cells.Item(79, 22).Value = "https://www.betexplorer.com/football/england/premier-league/west-ham-newcastle-utd/27P2HIIm/"

# Row 80: data now matches former row 79
$ws.Cells.Item(80, 6).Value = "Wolves"
$ws.Cells.Item(80, 7).Value = 1
$ws.Cells.Item(80, 8).Value = "Aston Villa"
$ws.Cells.Item(80, 9).Value = 1
$ws.Cells.Item(80, 10).Value = 2.81
$ws.Cells.Item(80, 11).Value = "24/09/2023 10:02"
$ws.Cells.Item(80, 12).Value = 3.7
$ws.Cells.Item(80, 13).Value = "08/10/2023 14:45"
$ws.Cells.Item(80, 14).Value = 3.39
$ws.Cells.Item(80, 15).Value = "24/09/2023 10:02"
$ws.Cells.Item(80, 16).Value = 3.68
$ws.Cells.Item(80, 17).Value = "08/10/2023 14:45"
$ws.Cells.Item(80, 18).Value = 2.62
$ws.Cells.Item(80, 19).Value = "24/09/2023 10:02"
$ws.Cells.Item(80, 20).Value = 2.08
$ws.Cells.Item(80, 21).Value = "08/10/2023 14:45"
$ws.Cells.Item(80, 22).Value = "https://www.betexplorer.com/football/england/premier-league/wolves-aston-villa/GAT6GxYg/"

# Row 83: data now matches former row 84
$ws.Cells.Item(83, 6).Value = "Nottingham"
$ws.Cells.Item(83, 7).Value = 2
$ws.Cells.Item(83, 8).Value = "Luton"
$ws.Cells.Item(83, 9).Value = 2
$ws.Cells.Item(83, 10).Value = 1.78
$ws.Cells.Item(83, 11).Value = "02/10/2023 08:30"
$ws.Cells.Item(83, 12).Value = 1.76
$ws.Cells.Item(83, 13).Value = "21/10/2023 15:56"
$ws.Cells.Item(83, 14).Value = 3.63
$ws.Cells.Item(83, 15).Value = "02/10/2023 08:30"
$ws.Cells.Item(83, 16).Value = 3.65
$ws.Cells.Item(83, 17).Value = "21/10/2023 15:58"
$ws.Cells.Item(83, 18).Value = 4.53
$ws.Cells.Item(83, 19).Value = "02/10/2023 08:30"
$ws.Cells.Item(83, 20).Value = 5.39
$ws.Cells.Item(83, 21).Value = "21/10/2023 15:58"
$ws.Cells.Item(83, 22).Value = "https://www.betexplorer.com/football/england/premier-league/nottingham-luton/tC3uVymm/"

# Row 84: data now matches former row 83
$ws.Cells.Item(84, 6).Value = "Newcastle"
$ws.Cells.Item(84, 7).Value = 4
$ws.Cells.Item(84, 8).Value = "Crystal Palace"
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 1.49
$ws.Cells.Item(84, 11).Value = "01/10/2023 23:01"
$ws.Cells.Item(84, 12).Value = 1.48
$ws.Cells.Item(84, 13).Value = "21/10/2023 15:50"
$ws.Cells.Item(84, 14).Value = 4.51
$ws.Cells.Item(84, 15).Value = "01/10/2023 23:01"
$ws.Cells.Item(84, 16).Value = 4.49
$ws.Cells.Item(84, 17).Value = "21/10/2023 15:58"
$ws.Cells.Item(84, 18).Value = 7.02
$ws.Cells.Item(84, 19).Value = "01/10/2023 23:01"
$ws.Cells.Item(84, 20).Value = 7.68
$ws.Cells.Item(84, 21).Value = "21/10/2023 15:58"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/england/premier-league/newcastle-utd-crystal-palace/2L4yWHXt/"

# Row 86: data now matches former row 87
$ws.Cells.Item(86, 6).Value = "Brentford"
$ws.Cells.Item(86, 7).Value = 3
$ws.Cells.Item(86, 8).Value = "Burnley"
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 1.63
$ws.Cells.Item(86, 11).Value = "02/10/2023 08:30"
$ws.Cells.Item(86, 12).Value = 1.74
$ws.Cells.Item(86, 13).Value = "21/10/2023 15:45"
$ws.Cells.Item(86, 14).Value = 3.97
$ws.Cells.Item(86, 15).Value = "02/10/2023 08:30"
$ws.Cells.Item(86, 16).Value = 3.9
$ws.Cells.Item(86, 17).Value = "21/10/2023 15:58"
$ws.Cells.Item(86, 18).Value = 5.07
$ws.Cells.Item(86, 19).Value = "02/10/2023 08:30"
$ws.Cells.Item(86, 20).Value = 5.09
$ws.Cells.Item(86, 21).Value = "21/10/2023 15:58"
$ws.Cells.Item(86, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brentford-burnley/6aMJDzIC/"

# Row 87: data now matches former row 86
$ws.Cells.Item(87, 6).Value = "Bournemouth"
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = "Wolves"
$ws.Cells.Item(87, 9).Value = 2
$ws.Cells.Item(87, 10).Value = 2.24
$ws.Cells.Item(87, 11).Value = "01/10/2023 23:01"
$ws.Cells.Item(87, 12).Value = 2.31
$ws.Cells.Item(87, 13).Value = "21/10/2023 15:55"
$ws.Cells.Item(87, 14).Value = 3.47
$ws.Cells.Item(87, 15).Value = "01/10/2023 23:01"
$ws.Cells.Item(87, 16).Value = 3.46
$ws.Cells.Item(87, 17).Value = "21/10/2023 15:57"
$ws.Cells.Item(87, 18).Value = 3.33
$ws.Cells.Item(87, 19).Value = "01/10/2023 23:01"
$ws.Cells.Item(87, 20).Value = 3.32
$ws.Cells.Item(87, 21).Value = "21/10/2023 15:58"
$ws.Cells.Item(87, 22).Value = "https://www.betexplorer.com/football/england/premier-league/bournemouth-wolves/bZIBFdm0/"

# Row 98: data now matches former row 99
$ws.Cells.Item(98, 6).Value = "Aston Villa"
$ws.Cells.Item(98, 7).Value = 3
$ws.Cells.Item(98, 8).Value = "Luton"
$ws.Cells.Item(98, 9).Value = 1
$ws.Cells.Item(98, 10).Value = 1.38
$ws.Cells.Item(98, 11).Value = "10/10/2023 14:32"
$ws.Cells.Item(98, 12).Value = 1.34
$ws.Cells.Item(98, 13).Value = "29/10/2023 14:55"
$ws.Cells.Item(98, 14).Value = 4.92
$ws.Cells.Item(98, 15).Value = "10/10/2023 14:32"
$ws.Cells.Item(98, 16).Value = 5.85
$ws.Cells.Item(98, 17).Value = "29/10/2023 14:55"
$ws.Cells.Item(98, 18).Value = 7.42
$ws.Cells.Item(98, 19).Value = "10/10/2023 14:32"
$ws.Cells.Item(98, 20).Value = 8.550000000000001
$ws.Cells.Item(98, 21).Value = "29/10/2023 14:57"
$ws.Cells.Item(98, 22).Value = "https://www.betexplorer.com/football/england/premier-league/aston-villa-luton/SSk1QD1I/"

# Row 99: data now matches former row 98
$ws.Cells.Item(99, 6).Value = "Brighton"
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = "Fulham"
$ws.Cells.Item(99, 9).Value = 1
$ws.Cells.Item(99, 10).Value = 1.51
$ws.Cells.Item(99, 11).Value = "10/10/2023 14:02"
$ws.Cells.Item(99, 12).Value = 1.64
$ws.Cells.Item(99, 13).Value = "29/10/2023 14:58"
$ws.Cells.Item(99, 14).Value = 4.73
$ws.Cells.Item(99, 15).Value = "10/10/2023 14:02"
$ws.Cells.Item(99, 16).Value = 4.32
$ws.Cells.Item(99, 17).Value = "29/10/2023 14:59"
$ws.Cells.Item(99, 18).Value = 6.2
$ws.Cells.Item(99, 19).Value = "10/10/2023 14:02"
$ws.Cells.Item(99, 20).Value = 5.29
$ws.Cells.Item(99, 21).Value = "29/10/2023 14:59"
$ws.Cells.Item(99, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brighton-fulham/6Jl5PXGO/"

# Row 103: data now matches former row 107
$ws.Cells.Item(103, 6).Value = "Sheffield Utd"
$ws.Cells.Item(103, 7).Value = 2
$ws.Cells.Item(103, 8).Value = "Wolves"
$ws.Cells.Item(103, 9).Value = 1
$ws.Cells.Item(103, 10).Value = 3.2
$ws.Cells.Item(103, 11).Value = "23/10/2023 15:48"
$ws.Cells.Item(103, 12).Value = 4.35
$ws.Cells.Item(103, 13).Value = "04/11/2023 15:59"
$ws.Cells.Item(103, 14).Value = 3.43
$ws.Cells.Item(103, 15).Value = "23/10/2023 15:48"
$ws.Cells.Item(103, 16).Value = 3.76
$ws.Cells.Item(103, 17).Value = "04/11/2023 15:59"
$ws.Cells.Item(103, 18).Value = 2.22
$ws.Cells.Item(103, 19).Value = "23/10/2023 15:48"
$ws.Cells.Item(103, 20).Value = 1.87
$ws.Cells.Item(103, 21).Value = "04/11/2023 15:59"
$ws.Cells.Item(103, 22).Value = "https://www.betexplorer.com/football/england/premier-league/sheffield-utd-wolves/0tW9gCV4/"

# Row 104: data now matches former row 106
$ws.Cells.Item(104, 6).Value = "Manchester City"
$ws.Cells.Item(104, 7).Value = 6
$ws.Cells.Item(104, 8).Value = "Bournemouth"
$ws.Cells.Item(104, 9).Value = 1
$ws.Cells.Item(104, 10).Value = 1.15
$ws.Cells.Item(104, 11).Value = "21/10/2023 22:01"
$ws.Cells.Item(104, 12).Value = 1.09
$ws.Cells.Item(104, 13).Value = "04/11/2023 15:50"
$ws.Cells.Item(104, 14).Value = 9.18
$ws.Cells.Item(104, 15).Value = "21/10/2023 22:01"
$ws.Cells.Item(104, 16).Value = 11.5
$ws.Cells.Item(104, 17).Value = "04/11/2023 15:21"
$ws.Cells.Item(104, 18).Value = 17.1
$ws.Cells.Item(104, 19).Value = "21/10/2023 22:01"
$ws.Cells.Item(104, 20).Value = 28.5
$ws.Cells.Item(104, 21).Value = "04/11/2023 15:52"
$ws.Cells.Item(104, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-city-bournemouth/AiwcdEon/"

# Row 105: data now matches former row 103
$ws.Cells.Item(105, 6).Value = "Brentford"
$ws.Cells.Item(105, 7).Value = 3
$ws.Cells.Item(105, 8).Value = "West Ham"
$ws.Cells.Item(105, 9).Value = 2
$ws.Cells.Item(105, 10).Value = 1.95
$ws.Cells.Item(105, 11).Value = "21/10/2023 20:02"
$ws.Cells.Item(105, 12).Value = 2.16
$ws.Cells.Item(105, 13).Value = "04/11/2023 15:50"
$ws.Cells.Item(105, 14).Value = 3.65
$ws.Cells.Item(105, 15).Value = "21/10/2023 20:02"
$ws.Cells.Item(105, 16).Value = 3.71
$ws.Cells.Item(105, 17).Value = "04/11/2023 15:50"
$ws.Cells.Item(105, 18).Value = 4
$ws.Cells.Item(105, 19).Value = "21/10/2023 20:02"
$ws.Cells.Item(105, 20).Value = 3.39
$ws.Cells.Item(105, 21).Value = "04/11/2023 15:50"
$ws.Cells.Item(105, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brentford-west-ham/MkBzuDGB/"

# Row 106: data now matches former row 104
$ws.Cells.Item(106, 6).Value = "Burnley"
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = "Crystal Palace"
$ws.Cells.Item(106, 9).Value = 2
$ws.Cells.Item(106, 10).Value = 2.77
$ws.Cells.Item(106, 11).Value = "23/10/2023 15:48"
$ws.Cells.Item(106, 12).Value = 3.53
$ws.Cells.Item(106, 13).Value = "04/11/2023 15:59"
$ws.Cells.Item(106, 14).Value = 3.22
$ws.Cells.Item(106, 15).Value = "23/10/2023 15:48"
$ws.Cells.Item(106, 16).Value = 3.21
$ws.Cells.Item(106, 17).Value = "04/11/2023 15:59"
$ws.Cells.Item(106, 18).Value = 2.61
$ws.Cells.Item(106, 19).Value = "23/10/2023 15:48"
$ws.Cells.Item(106, 20).Value = 2.31
$ws.Cells.Item(106, 21).Value = "04/11/2023 15:59"
$ws.Cells.Item(106, 22).Value = "https://www.betexplorer.com/football/england/premier-league/burnley-crystal-palace/0tAvvXVH/"

# Row 107: data now matches former row 105
$ws.Cells.Item(107, 6).Value = "Everton"
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = "Brighton"
$ws.Cells.Item(107, 9).Value = 1
$ws.Cells.Item(107, 10).Value = 2.86
$ws.Cells.Item(107, 11).Value = "21/10/2023 20:02"
$ws.Cells.Item(107, 12).Value = 2.86
$ws.Cells.Item(107, 13).Value = "04/11/2023 15:59"
$ws.Cells.Item(107, 14).Value = 3.84
$ws.Cells.Item(107, 15).Value = "21/10/2023 20:02"
$ws.Cells.Item(107, 16).Value = 3.52
$ws.Cells.Item(107, 17).Value = "04/11/2023 15:59"
$ws.Cells.Item(107, 18).Value = 2.24
$ws.Cells.Item(107, 19).Value = "21/10/2023 20:02"
$ws.Cells.Item(107, 20).Value = 2.54
$ws.Cells.Item(107, 21).Value = "04/11/2023 15:59"
$ws.Cells.Item(107, 22).Value = "https://www.betexplorer.com/football/england/premier-league/everton-brighton/f39rwioO/"

# Row 113: data now matches former row 115
$ws.Cells.Item(113, 6).Value = "Arsenal"
$ws.Cells.Item(113, 7).Value = 3
$ws.Cells.Item(113, 8).Value = "Burnley"
$ws.Cells.Item(113, 9).Value = 1
$ws.Cells.Item(113, 10).Value = 1.24
$ws.Cells.Item(113, 11).Value = "29/10/2023 11:22"
$ws.Cells.Item(113, 12).Value = 1.19
$ws.Cells.Item(113, 13).Value = "11/11/2023 15:13"
$ws.Cells.Item(113, 14).Value = 6.5
$ws.Cells.Item(113, 15).Value = "29/10/2023 11:22"
$ws.Cells.Item(113, 16).Value = 7.49
$ws.Cells.Item(113, 17).Value = "11/11/2023 15:36"
$ws.Cells.Item(113, 18).Value = 12.06
$ws.Cells.Item(113, 19).Value = "29/10/2023 11:22"
$ws.Cells.Item(113, 20).Value = 16.59
$ws.Cells.Item(113, 21).Value = "11/11/2023 15:36"
$ws.Cells.Item(113, 22).Value = "https://www.betexplorer.com/football/england/premier-league/arsenal-burnley/ncYLjAFN/"

# Row 115: data now matches former row 113
$ws.Cells.Item(115, 6).Value = "Manchester Utd"
$ws.Cells.Item(115, 7).Value = 1
$ws.Cells.Item(115, 8).Value = "Luton"
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 10).Value = 1.2
$ws.Cells.Item(115, 11).Value = "29/10/2023 11:21"
$ws.Cells.Item(115, 12).Value = 1.29
$ws.Cells.Item(115, 13).Value = "11/11/2023 15:58"
$ws.Cells.Item(115, 14).Value = 7.39
$ws.Cells.Item(115, 15).Value = "29/10/2023 11:21"
$ws.Cells.Item(115, 16).Value = 6.14
$ws.Cells.Item(115, 17).Value = "11/11/2023 15:59"
$ws.Cells.Item(115, 18).Value = 13.36
$ws.Cells.Item(115, 19).Value = "29/10/2023 11:21"
$ws.Cells.Item(115, 20).Value = 10.06
$ws.Cells.Item(115, 21).Value = "11/11/2023 15:59"
$ws.Cells.Item(115, 22).Value = "https://www.betexplorer.com/football/england/premier-league/manchester-united-luton/4lXv7Va4/"

# Row 117: data now matches former row 119
$ws.Cells.Item(117, 6).Value = "Liverpool"
$ws.Cells.Item(117, 7).Value = 3
$ws.Cells.Item(117, 8).Value = "Brentford"
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 1.48
$ws.Cells.Item(117, 11).Value = "28/10/2023 22:02"
$ws.Cells.Item(117, 12).Value = 1.4
$ws.Cells.Item(117, 13).Value = "12/11/2023 14:56"
$ws.Cells.Item(117, 14).Value = 4.87
$ws.Cells.Item(117, 15).Value = "28/10/2023 22:02"
$ws.Cells.Item(117, 16).Value = 5.51
$ws.Cells.Item(117, 17).Value = "12/11/2023 14:56"
$ws.Cells.Item(117, 18).Value = 6.16
$ws.Cells.Item(117, 19).Value = "28/10/2023 22:02"
$ws.Cells.Item(117, 20).Value = 7.48
$ws.Cells.Item(117, 21).Value = "12/11/2023 14:59"
$ws.Cells.Item(117, 22).Value = "https://www.betexplorer.com/football/england/premier-league/liverpool-brentford/QcYz8Bpb/"

# Row 118: data now matches former row 120
$ws.Cells.Item(118, 6).Value = "West Ham"
$ws.Cells.Item(118, 7).Value = 3
$ws.Cells.Item(118, 8).Value = "Nottingham"
$ws.Cells.Item(118, 9).Value = 2
$ws.Cells.Item(118, 10).Value = 1.82
$ws.Cells.Item(118, 11).Value = "28/10/2023 22:02"
$ws.Cells.Item(118, 12).Value = 1.81
$ws.Cells.Item(118, 13).Value = "12/11/2023 14:35"
$ws.Cells.Item(118, 14).Value = 3.74
$ws.Cells.Item(118, 15).Value = "28/10/2023 22:02"
$ws.Cells.Item(118, 16).Value = 3.84
$ws.Cells.Item(118, 17).Value = "12/11/2023 14:35"
$ws.Cells.Item(118, 18).Value = 4.46
$ws.Cells.Item(118, 19).Value = "28/10/2023 22:02"
$ws.Cells.Item(118, 20).Value = 4.59
$ws.Cells.Item(118, 21).Value = "12/11/2023 14:58"
$ws.Cells.Item(118, 22).Value = "https://www.betexplorer.com/football/england/premier-league/west-ham-nottingham/YNyq6kFA/"

# Row 119: data now matches former row 117
$ws.Cells.Item(119, 6).Value = "Aston Villa"
$ws.Cells.Item(119, 7).Value = 3
$ws.Cells.Item(119, 8).Value = "Fulham"
$ws.Cells.Item(119, 9).Value = 1
$ws.Cells.Item(119, 10).Value = 1.59
$ws.Cells.Item(119, 11).Value = "29/10/2023 00:02"
$ws.Cells.Item(119, 12).Value = 1.64
$ws.Cells.Item(119, 13).Value = "12/11/2023 14:58"
$ws.Cells.Item(119, 14).Value = 4.16
$ws.Cells.Item(119, 15).Value = "29/10/2023 00:02"
$ws.Cells.Item(119, 16).Value = 4.22
$ws.Cells.Item(119, 17).Value = "12/11/2023 14:55"
$ws.Cells.Item(119, 18).Value = 5.06
$ws.Cells.Item(119, 19).Value = "29/10/2023 00:02"
$ws.Cells.Item(119, 20).Value = 5.4
$ws.Cells.Item(119, 21).Value = "12/11/2023 14:58"
$ws.Cells.Item(119, 22).Value = "https://www.betexplorer.com/football/england/premier-league/aston-villa-fulham/EJ4IkUUT/"

# Row 120: data now matches former row 118
$ws.Cells.Item(120, 6).Value = "Brighton"
$ws.Cells.Item(120, 7).Value = 1
$ws.Cells.Item(120, 8).Value = "Sheffield Utd"
$ws.Cells.Item(120, 9).Value = 1
$ws.Cells.Item(120, 10).Value = 1.23
$ws.Cells.Item(120, 11).Value = "29/10/2023 11:22"
$ws.Cells.Item(120, 12).Value = 1.27
$ws.Cells.Item(120, 13).Value = "12/11/2023 14:55"
$ws.Cells.Item(120, 14).Value = 6.85
$ws.Cells.Item(120, 15).Value = "29/10/2023 11:22"
$ws.Cells.Item(120, 16).Value = 6.23
$ws.Cells.Item(120, 17).Value = "12/11/2023 14:59"
$ws.Cells.Item(120, 18).Value = 11.63
$ws.Cells.Item(120, 19).Value = "29/10/2023 11:22"
$ws.Cells.Item(120, 20).Value = 11.19
$ws.Cells.Item(120, 21).Value = "12/11/2023 14:59"
$ws.Cells.Item(120, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brighton-sheffield-utd/GSjpolwo/"

# Row 123: data now matches former row 127
$ws.Cells.Item(123, 6).Value = "Burnley"
$ws.Cells.Item(123, 7).Value = 1
$ws.Cells.Item(123, 8).Value = "West Ham"
$ws.Cells.Item(123, 9).Value = 2
$ws.Cells.Item(123, 10).Value = 2.61
$ws.Cells.Item(123, 11).Value = "05/11/2023 11:03"
$ws.Cells.Item(123, 12).Value = 3.48
$ws.Cells.Item(123, 13).Value = "25/11/2023 15:56"
$ws.Cells.Item(123, 14).Value = 3.38
$ws.Cells.Item(123, 15).Value = "05/11/2023 11:03"
$ws.Cells.Item(123, 16).Value = 3.43
$ws.Cells.Item(123, 17).Value = "25/11/2023 15:59"
$ws.Cells.Item(123, 18).Value = 2.78
$ws.Cells.Item(123, 19).Value = "05/11/2023 11:03"
$ws.Cells.Item(123, 20).Value = 2.23
$ws.Cells.Item(123, 21).Value = "25/11/2023 15:56"
$ws.Cells.Item(123, 22).Value = "https://www.betexplorer.com/football/england/premier-league/burnley-west-ham/4nhLn1Ek/"

# Row 124: data now matches former row 125
$ws.Cells.Item(124, 6).Value = "Luton"
$ws.Cells.Item(124, 7).Value = 2
$ws.Cells.Item(124, 8).Value = "Crystal Palace"
$ws.Cells.Item(124, 9).Value = 1
$ws.Cells.Item(124, 10).Value = 3.35
$ws.Cells.Item(124, 11).Value = "05/11/2023 11:03"
$ws.Cells.Item(124, 12).Value = 4.12
$ws.Cells.Item(124, 13).Value = "25/11/2023 15:57"
$ws.Cells.Item(124, 14).Value = 3.27
$ws.Cells.Item(124, 15).Value = "05/11/2023 11:03"
$ws.Cells.Item(124, 16).Value = 3.42
$ws.Cells.Item(124, 17).Value = "25/11/2023 15:59"
$ws.Cells.Item(124, 18).Value = 2.3
$ws.Cells.Item(124, 19).Value = "05/11/2023 11:03"
$ws.Cells.Item(124, 20).Value = 2.03
$ws.Cells.Item(124, 21).Value = "25/11/2023 15:59"
$ws.Cells.Item(124, 22).Value = "https://www.betexplorer.com/football/england/premier-league/luton-crystal-palace/EHZB9OL2/"

# Row 125: data now matches former row 124
$ws.Cells.Item(125, 6).Value = "Newcastle"
$ws.Cells.Item(125, 7).Value = 4
$ws.Cells.Item(125, 8).Value = "Chelsea"
$ws.Cells.Item(125, 9).Value = 1
$ws.Cells.Item(125, 10).Value = 1.84
$ws.Cells.Item(125, 11).Value = "05/11/2023 11:03"
$ws.Cells.Item(125, 12).Value = 2.6
$ws.Cells.Item(125, 13).Value = "25/11/2023 15:25"
$ws.Cells.Item(125, 14).Value = 3.75
$ws.Cells.Item(125, 15).Value = "05/11/2023 11:03"
$ws.Cells.Item(125, 16).Value = 3.55
$ws.Cells.Item(125, 17).Value = "25/11/2023 15:07"
$ws.Cells.Item(125, 18).Value = 4.33
$ws.Cells.Item(125, 19).Value = "05/11/2023 11:03"
$ws.Cells.Item(125, 20).Value = 2.76
$ws.Cells.Item(125, 21).Value = "25/11/2023 15:25"
$ws.Cells.Item(125, 22).Value = "https://www.betexplorer.com/football/england/premier-league/newcastle-utd-chelsea/nVXJ72jF/"

# Row 126: data now matches former row 123
$ws.Cells.Item(126, 6).Value = "Nottingham"
$ws.Cells.Item(126, 7).Value = 2
$ws.Cells.Item(126, 8).Value = "Brighton"
$ws.Cells.Item(126, 9).Value = 3
$ws.Cells.Item(126, 10).Value = 3.29
$ws.Cells.Item(126, 11).Value = "05/11/2023 11:03"
$ws.Cells.Item(126, 12).Value = 3.69
$ws.Cells.Item(126, 13).Value = "25/11/2023 15:57"
$ws.Cells.Item(126, 14).Value = 3.61
$ws.Cells.Item(126, 15).Value = "05/11/2023 11:03"
$ws.Cells.Item(126, 16).Value = 3.63
$ws.Cells.Item(126, 17).Value = "25/11/2023 15:57"
$ws.Cells.Item(126, 18).Value = 2.18
$ws.Cells.Item(126, 19).Value = "05/11/2023 11:03"
$ws.Cells.Item(126, 20).Value = 2.07
$ws.Cells.Item(126, 21).Value = "25/11/2023 15:57"
$ws.Cells.Item(126, 22).Value = "https://www.betexplorer.com/football/england/premier-league/nottingham-brighton/0YyO6M5L/"

# Row 127: data now matches former row 126
$ws.Cells.Item(127, 6).Value = "Sheffield Utd"
$ws.Cells.Item(127, 7).Value = 1
$ws.Cells.Item(127, 8).Value = "Bournemouth"
$ws.Cells.Item(127, 9).Value = 3
$ws.Cells.Item(127, 10).Value = 2.68
$ws.Cells.Item(127, 11).Value = "05/11/2023 11:03"
$ws.Cells.Item(127, 12).Value = 3.57
$ws.Cells.Item(127, 13).Value = "25/11/2023 15:52"
$ws.Cells.Item(127, 14).Value = 3.4
$ws.Cells.Item(127, 15).Value = "05/11/2023 11:03"
$ws.Cells.Item(127, 16).Value = 3.5
$ws.Cells.Item(127, 17).Value = "25/11/2023 15:52"
$ws.Cells.Item(127, 18).Value = 2.68
$ws.Cells.Item(127, 19).Value = "05/11/2023 11:03"
$ws.Cells.Item(127, 20).Value = 2.16
$ws.Cells.Item(127, 21).Value = "25/11/2023 15:52"
$ws.Cells.Item(127, 22).Value = "https://www.betexplorer.com/football/england/premier-league/sheffield-utd-bournemouth/rcRS5tLR/"

# --- Part 2: append 4 new match rows (132-135), copying style from row 131 for columns A and E ---

# Row 132
$ws.Cells.Item(131, 1).Copy($ws.Cells.Item(132, 1))
$ws.Cells.Item(131, 5).Copy($ws.Cells.Item(132, 5))
$ws.Cells.Item(132, 1).Value = 131
$ws.Cells.Item(132, 2).Value = "england"
$ws.Cells.Item(132, 3).Value = "premier-league"
$ws.Cells.Item(132, 4).Value = "2023-2024"
$ws.Cells.Item(132, 5).Value = 45262.66666666666
$ws.Cells.Item(132, 6).Value = "Burnley"
$ws.Cells.Item(132, 7).Value = 5
$ws.Cells.Item(132, 8).Value = "Sheffield Utd"
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 1.7
$ws.Cells.Item(132, 11).Value = "15/11/2023 16:01"
$ws.Cells.Item(132, 12).Value = 1.74
$ws.Cells.Item(132, 13).Value = "02/12/2023 15:54"
$ws.Cells.Item(132, 14).Value = 3.93
$ws.Cells.Item(132, 15).Value = "15/11/2023 16:01"
$ws.Cells.Item(132, 16).Value = 3.81
$ws.Cells.Item(132, 17).Value = "02/12/2023 15:58"
$ws.Cells.Item(132, 18).Value = 5.01
$ws.Cells.Item(132, 19).Value = "15/11/2023 16:01"
$ws.Cells.Item(132, 20).Value = 5.19
$ws.Cells.Item(132, 21).Value = "02/12/2023 15:58"
$ws.Cells.Item(132, 22).Value = "https://www.betexplorer.com/football/england/premier-league/burnley-sheffield-utd/Sl5BcIS7/"

# Row 133
$ws.Cells.Item(131, 1).Copy($ws.Cells.Item(133, 1))
$ws.Cells.Item(131, 5).Copy($ws.Cells.Item(133, 5))
$ws.Cells.Item(133, 1).Value = 132
$ws.Cells.Item(133, 2).Value = "england"
$ws.Cells.Item(133, 3).Value = "premier-league"
$ws.Cells.Item(133, 4).Value = "2023-2024"
$ws.Cells.Item(133, 5).Value = 45262.66666666666
$ws.Cells.Item(133, 6).Value = "Arsenal"
$ws.Cells.Item(133, 7).Value = 2
$ws.Cells.Item(133, 8).Value = "Wolves"
$ws.Cells.Item(133, 9).Value = 1
$ws.Cells.Item(133, 10).Value = 1.27
$ws.Cells.Item(133, 11).Value = "15/11/2023 16:01"
$ws.Cells.Item(133, 12).Value = 1.24
$ws.Cells.Item(133, 13).Value = "02/12/2023 15:08"
$ws.Cells.Item(133, 14).Value = 5.64
$ws.Cells.Item(133, 15).Value = "15/11/2023 16:01"
$ws.Cells.Item(133, 16).Value = 6.5
$ws.Cells.Item(133, 17).Value = "02/12/2023 15:10"
$ws.Cells.Item(133, 18).Value = 9.24
$ws.Cells.Item(133, 19).Value = "15/11/2023 16:01"
$ws.Cells.Item(133, 20).Value = 13
$ws.Cells.Item(133, 21).Value = "02/12/2023 15:11"
$ws.Cells.Item(133, 22).Value = "https://www.betexplorer.com/football/england/premier-league/arsenal-wolves/WMB2avbe/"

# Row 134
$ws.Cells.Item(131, 1).Copy($ws.Cells.Item(134, 1))
$ws.Cells.Item(131, 5).Copy($ws.Cells.Item(134, 5))
$ws.Cells.Item(134, 1).Value = 133
$ws.Cells.Item(134, 2).Value = "england"
$ws.Cells.Item(134, 3).Value = "premier-league"
$ws.Cells.Item(134, 4).Value = "2023-2024"
$ws.Cells.Item(134, 5).Value = 45262.66666666666
$ws.Cells.Item(134, 6).Value = "Brentford"
$ws.Cells.Item(134, 7).Value = 3
$ws.Cells.Item(134, 8).Value = "Luton"
$ws.Cells.Item(134, 9).Value = 1
$ws.Cells.Item(134, 10).Value = 1.42
$ws.Cells.Item(134, 11).Value = "15/11/2023 16:01"
$ws.Cells.Item(134, 12).Value = 1.54
$ws.Cells.Item(134, 13).Value = "02/12/2023 15:33"
$ws.Cells.Item(134, 14).Value = 4.64
$ws.Cells.Item(134, 15).Value = "15/11/2023 16:01"
$ws.Cells.Item(134, 16).Value = 4.43
$ws.Cells.Item(134, 17).Value = "02/12/2023 15:58"
$ws.Cells.Item(134, 18).Value = 6.76
$ws.Cells.Item(134, 19).Value = "15/11/2023 16:01"
$ws.Cells.Item(134, 20).Value = 6.38
$ws.Cells.Item(134, 21).Value = "02/12/2023 15:58"
$ws.Cells.Item(134, 22).Value = "https://www.betexplorer.com/football/england/premier-league/brentford-luton/Eo17bbD1/"

# Row 135
$ws.Cells.Item(131, 1).Copy($ws.Cells.Item(135, 1))
$ws.Cells.Item(131, 5).Copy($ws.Cells.Item(135, 5))
$ws.Cells.Item(135, 1).Value = 134
$ws.Cells.Item(135, 2).Value = "england"
$ws.Cells.Item(135, 3).Value = "premier-league"
$ws.Cells.Item(135, 4).Value = "2023-2024"
$ws.Cells.Item(135, 5).Value = 45262.77083333334
$ws.Cells.Item(135, 6).Value = "Nottingham"
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = "Everton"
$ws.Cells.Item(135, 9).Value = 1
$ws.Cells.Item(135, 10).Value = 2.67
$ws.Cells.Item(135, 11).Value = "15/11/2023 16:01"
$ws.Cells.Item(135, 12).Value = 2.88
$ws.Cells.Item(135, 13).Value = "02/12/2023 18:26"
$ws.Cells.Item(135, 14).Value = 3.28
$ws.Cells.Item(135, 15).Value = "15/11/2023 16:01"
$ws.Cells.Item(135, 16).Value = 3.16
$ws.Cells.Item(135, 17).Value = "02/12/2023 18:27"
$ws.Cells.Item(135, 18).Value = 2.63
$ws.Cells.Item(135, 19).Value = "15/11/2023 16:01"
$ws.Cells.Item(135, 20).Value = 2.76
$ws.Cells.Item(135, 21).Value = "02/12/2023 17:39"
$ws.Cells.Item(135, 22).Value = "https://www.betexplorer.com/football/england/premier-league/nottingham-everton/YkYXEJcr/"

# --- Part 3: update the sheet dimension reference to cover the new rows ---
$ws.Range("A1:V135").Select() | Out-Null